$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear C2 entirely (the value is removed in the corrected data)
$ws.Range("C2").Value = $null

# Updated y_1_forecast (column E) values - tiny precision corrections from bug fix
$ws.Range("E2").Value = 2.228397109637203
$ws.Range("E3").Value = -0.2018858887078645
$ws.Range("E4").Value = 0.2108047537406454
$ws.Range("E5").Value = 0.96842791562195
$ws.Range("E6").Value = 1.698182372097512
$ws.Range("E8").Value = 1.325176859452348
$ws.Range("E9").Value = 1.522808462763692
$ws.Range("E10").Value = 1.693557061600948
$ws.Range("E11").Value = 1.905564797014669
$ws.Range("E12").Value = 2.08247707460909
$ws.Range("E13").Value = 1.845103901518907
$ws.Range("E14").Value = 0.02570757229445331

# Updated y_0_forecast (column C) values - tiny precision corrections from bug fix
$ws.Range("C3").Value = -0.4626567965509643
$ws.Range("C7").Value = -0.2674335569108788
$ws.Range("C8").Value = 2.038609866767915
$ws.Range("C10").Value = 2.246337373618967
$ws.Range("C11").Value = 1.777150434343544
$ws.Range("C12").Value = 1.741137453897323
$ws.Range("C13").Value = 1.562095320687429
$ws.Range("C17").Value = -0.2883789941992232
$ws.Range("C19").Value = 1.431852292002245
